# Add new metadata report row for Akurana (update row 2 values/labels)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A2").Value = 2025
$ws.Range("B2").Value = "JAN"
$ws.Range("C2").Value = "10/01-01/01"
$ws.Range("D2").Value = "11/01, 12/01, 13/01, 14/01, 15/01, 16/01, 17/01, 18/01, 19/01, 20/01, 21/01, 22/01, 23/01, 24/01, 25/01, 26/01, 27/01, 28/01, 29/01, 30/01, 31/01"

$ws.Range("F2").Value = 22.83820662768031
$ws.Range("G2").Value = 23.61477582846004
$ws.Range("H2").Value = 20.69738791423002
$ws.Range("I2").Value = 0.3512670565302144
$ws.Range("J2").Value = 1.196296296296296
$ws.Range("K2").Value = 6.302690058479532
$ws.Range("L2").Value = 118.8674463937622
$ws.Range("M2").Value = 3.072358674463938
$ws.Range("N2").Value = 1.381247563352827
$ws.Range("O2").Value = 19.78253411306043
$ws.Range("P2").Value = 27.17270955165692
$ws.Range("Q2").Value = 27.17270955165692
$ws.Range("R2").Value = 961.310760233918
$ws.Range("S2").Value = 88.70409356725146
$ws.Range("T2").Value = 0.6923976608187135
$ws.Range("U2").Value = 91.63009746588695
$ws.Range("V2").Value = 23.20931773879142
$ws.Range("W2").Value = 85.32007797270956
$ws.Range("Y2").Value = 23.80615984405458
$ws.Range("Z2").Value = 20.58249512670565
$ws.Range("AA2").Value = 961.310760233918
$ws.Range("AC2").Value = 116.729044834308
$ws.Range("AD2").Value = 0.2131384015594542
